$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Cape Verde - Eswatini -> final score 3:0, checkmark result
$ws.Range("B2").Value = "Cape Verde ✓ - Eswatini: 3:0"
$ws.Range("G2").Value = "✓"
$ws.Range("H2").Value = 3
# I2 stays FALSE (unchanged)

# Row 3: Iceland - France -> final score 2:2
$ws.Range("B3").Value = "Iceland - France : 2:2"
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = $true

# Row 4: Cameroon - Angola -> final score 0:0
$ws.Range("B4").Value = "Cameroon  - Angola: 0:0"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = $true

# Row 5: North Macedonia - Kazakhstan -> final score 1:1
$ws.Range("B5").Value = "North Macedonia  - Kazakhstan: 1:1"
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = $true

# Row 6: Sweden - Kosovo -> final score 0:1, X result marker
$ws.Range("B6").Value = "Sweden X - Kosovo: 0:1"
$ws.Range("G6").Value = "X"
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = $true

# Row 7: The Strongest La Paz - Club Jorge Wilstermann (just extra space in label)
$ws.Range("B7").Value = "The Strongest La Paz  - Club Jorge Wilstermann: 1:1"

# Row 8: Tunisia - Namibia -> final score 3:0, checkmark result
$ws.Range("B8").Value = "Tunisia ✓ - Namibia: 3:0"
$ws.Range("G8").Value = "✓"
$ws.Range("H8").Value = 3
# I8 stays FALSE (unchanged)

# Row 9 (Club Universidad de Chile fixture) is removed entirely; the
# summary formulas below shift up one row and their references adjust
# automatically.
$ws.Rows("9").Delete()
